# Update the Metadata sheet: Date and Contact values changed.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B8").Value = "2025-05-12T20:06:16+02:00"
$meta.Range("B10").Value = "Constantin Renner (http://example.org/example-publisher, constantinrenner1@gmail.com)"

# Update the Include #0 sheet: insert a new concept row before the
# trailing blank separator / "System URI" rows (which both shift down
# by one row).
$inc = $wb.Worksheets.Item("Include #0")

$inc.Rows.Item(12).Insert()

# Write the new concept's values. The leading apostrophe forces the
# numeric-looking code to be stored as text (matching how the other
# SNOMED codes in this column are stored), and is stripped from the
# actual cell content.
$inc.Range("A12").Value = "'119567009"
$inc.Range("B12").Value = "Structure of artery of head (body structure)"

# Re-apply the row formatting (border/alignment) from the row below so
# the new row matches the existing table styling instead of picking up
# a generic inserted-row format.
$inc.Range("A13:B13").Copy()
$inc.Range("A12:B12").PasteSpecial(-4122)
